# Change the computations of the KPIs
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update StartingInventories (C) and SetupCosts (E) ---
$wsProd = $wb.Worksheets.Item("Productdata")

$wsProd.Range("E2").Value = 0.1111111111111111

$wsProd.Range("C3").Value = 0
$wsProd.Range("E3").Value = 0.02570833333333333

$wsProd.Range("C4").Value = 0
$wsProd.Range("E4").Value = 0.02569444444444444

$wsProd.Range("C5").Value = 0
$wsProd.Range("E5").Value = 0.02565277777777777

$wsProd.Range("C6").Value = 0
$wsProd.Range("E6").Value = 0.05130555555555555

$wsProd.Range("E7").Value = 0.02570833333333333

$wsProd.Range("E8").Value = 0.02569444444444444

$wsProd.Range("E9").Value = 0.05130555555555555

$wsProd.Range("E10").Value = 0.02565277777777777

$wsProd.Range("C11").Value = 0
$wsProd.Range("E11").Value = 0.022375

$wsProd.Range("C12").Value = 0
$wsProd.Range("E12").Value = 0.02236111111111111

$wsProd.Range("C13").Value = 0
$wsProd.Range("E13").Value = 0.02231944444444444

$wsProd.Range("C14").Value = 0
$wsProd.Range("E14").Value = 0.04463888888888888

# --- Capacity sheet: update capacity values (B) ---
$wsCap = $wb.Worksheets.Item("Capacity")

$wsCap.Range("B2").Value = 125
$wsCap.Range("B3").Value = 5
$wsCap.Range("B4").Value = 10
$wsCap.Range("B5").Value = 25
$wsCap.Range("B6").Value = 40
$wsCap.Range("B7").Value = 10
$wsCap.Range("B8").Value = 15
$wsCap.Range("B9").Value = 15
$wsCap.Range("B10").Value = 25
$wsCap.Range("B11").Value = 25
$wsCap.Range("B12").Value = 20
$wsCap.Range("B13").Value = 25
$wsCap.Range("B14").Value = 20

# --- ProcessingTime sheet: update processing time values ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$wsProc.Range("B2").Value = 5
$wsProc.Range("D4").Value = 2
$wsProc.Range("E5").Value = 5
$wsProc.Range("J10").Value = 5
$wsProc.Range("K11").Value = 5
$wsProc.Range("M13").Value = 5
$wsProc.Range("N14").Value = 2
